$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1250
$ws.Range("C2").Value = 6000

$ws.Range("B3").Value = 2500
$ws.Range("C3").Value = 1250

$ws.Range("B4").Value = 1250
$ws.Range("C4").Value = 6000

$ws.Range("B5").Value = 2500
$ws.Range("C5").Value = 1250

$ws.Range("B6").Value = 2300
$ws.Range("C6").Value = 6000

$ws.Range("B7").Value = 2300
$ws.Range("C7").Value = 1150

$ws.Range("B8").Value = 2300
$ws.Range("C8").Value = 1150
